# Fix bug in logCommonResults() method inside SharedMethods.java class
#
# Mirrors the authored change to testData.xlsx:
#   - Google sheet: move the saved cursor/selection to D7 (sheet stays inactive)
#   - DuckDuckGo sheet: move the saved cursor/selection to C9, widen column A,
#     and replace the "hire me please" test-keyword row with the existing
#     "ruby is the best programming language" keyword (the now-unused shared
#     string drops out of the workbook automatically)

$wb = $excel.ActiveWorkbook

$google = $wb.Worksheets.Item("Google")
$duck   = $wb.Worksheets.Item("DuckDuckGo")

# Update the (inactive) Google sheet's remembered selection without
# stealing the active-tab flag from DuckDuckGo.
$null = $google.Range("D7").Select()

# Re-activate DuckDuckGo (it was, and should remain, the active tab) and
# apply its changes.
$null = $duck.Activate()

# A4 held the leftover "hire me please" keyword - replace it with the
# legitimate keyword value already used elsewhere in the sheet.
$duck.Range("A4").Value = "ruby is the best programming language"

# Widen column A to (the closest this host's pixel-quantized ColumnWidth
# setter can represent to) 36.85546875 characters.
$duck.Columns.Item(1).ColumnWidth = 36

# Finally, move DuckDuckGo's own saved selection to C9.
$null = $duck.Range("C9").Select()
